$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, bordered, centered) from H1 onto the two
# new header cells so I1/J1 match the rest of row 1's formatting.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data columns I and J for rows 2-10
$data = @{
    2  = @(3, 8)
    3  = @(5, 7)
    4  = @(1, 8)
    5  = @(1, 7)
    6  = @(1, 6)
    7  = @(1, 6)
    8  = @(6, 7)
    9  = @(8, 9)
    10 = @(3, 4)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
